# Update cryptos list with latest scraped prices / 1h volume changes.
# Generated for commit: "Updated cryptos list on Fri Mar  8 13:45:27 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.499.92"
$ws.Range("E2").Value = "  +1.03%  "

# Row 3
$ws.Range("D3").Value = "3.974.27"
$ws.Range("E3").Value = "  +4.93%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "485.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.44%  "

# Row 7
$ws.Range("E7").Value = "  +1.19%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.08%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.733"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.42%  "

# Row 10
$ws.Range("E10").Value = "  +12.54%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000356"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +15.79%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.29%  "

# Row 13
$ws.Range("D13").Value = "4.602.68"
$ws.Range("E13").Value = "  +4.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.56%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.76%  "

# Row 16
$ws.Range("D16").Value = "3.994.62"
$ws.Range("E16").Value = "  +5.37%  "

# Row 17
$ws.Range("E17").Value = "  +0.52%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.34%  "

# Row 19
$ws.Range("E19").Value = "  +1.77%  "

# Row 20
$ws.Range("D20").Value = "67.683.62"
$ws.Range("E20").Value = "  +1.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "435.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.30%  "

# Row 22
$ws.Range("E22").Value = "  +5.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.24%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.09%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.40%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.81%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "729.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "

# Row 30
$ws.Range("E30").Value = "  -1.16%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.44"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.79%  "

# Row 32
$ws.Range("E32").Value = "  +4.55%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "42.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.38%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0857"
$ws.Range("E34").Value = "  +28.20%  "

# Row 35
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.28%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.154"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.99%  "

# Row 37
$ws.Range("E37").Value = "  -1.02%  "

# Row 38
$ws.Range("E38").Value = "  -0.07%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0477"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.26%  "

# Row 41
$ws.Range("E41").Value = "  +2.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.340"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.09%  "

# Row 45
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.88%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.01%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.59%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "149.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.55%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.61%  "
